# Add a new, empty placeholder paragraph right after the existing
# "{/boletas}" closing-tag paragraph (the last paragraph in the body),
# immediately before the section properties.  The new paragraph inherits
# the same indentation / centering / run-size formatting as the paragraph
# it follows.

$d = $word.ActiveDocument

$lastParagraph = $d.Paragraphs.Last
$lastParagraph.Range.InsertParagraphAfter()
